$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.525.82"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "2.347.74"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.32"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.44"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.04%  "
$ws.Range("E7").Value = "  -3.43%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.508"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.25"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.81%  "
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("E13").Value = "  -2.90%  "
$ws.Range("D14").Value = "2.703.48"
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.76"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "2.313.42"
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.810"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "43.443.88"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.84"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.12"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.28"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.56"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.02%  "
$ws.Range("E24").Value = "  -3.59%  "
$ws.Range("E25").Value = "  -3.45%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.10"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.31%  "
$ws.Range("E28").Value = "  +3.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.74"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "165.80"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.26"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  -4.14%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.42"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.76%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.51"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("E36").Value = "  -8.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0708"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.93"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -6.94%  "
$ws.Range("E39").Value = "  -6.94%  "
$ws.Range("E40").Value = "  -3.43%  "
$ws.Range("E41").Value = "  -3.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.45"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("D43").Value = "1.983.72"
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.72"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -8.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.06"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.95"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.52"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.63%  "
$ws.Range("E49").Value = "  +3.04%  "
$ws.Range("D50").Value = "2.568.03"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("E51").Value = "  -1.83%  "
